$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$accent6 = [Microsoft.Office.Core.MsoThemeColorIndex]::msoThemeColorAccent6

# --- Part 1: English comment #6 (shape "내용 개체 틀 8") ---
$shEn = $s.Shapes.Item(2)
$trEn = $shEn.TextFrame.TextRange

$oldEn = "6. The re-ranking performance should explain the need for Lattice + Transformer (En)? and why transformers are used? How performance is increased using Lattice + Transformer (En). What does the cost of performance, mean with the same computational resources? The Lattice + Transformer (En) model needs to be explained more in detail with diagrams."
$newEn = "6. The re-ranking performance should explain the need for Lattice + Transformer (Encoder)? and why transformers are used? How performance is increased using Lattice + Transformer (Encoder). What does the cost of performance, mean with the same computational resources? The Lattice + Transformer (Encoder) model needs to be explained more in detail with diagrams."

$foundEn = $trEn.Find($oldEn)
if ($foundEn -eq $null) {
  throw "Could not find the English comment #6 paragraph text"
}
$foundEn.Text = $newEn

# Re-find (text length changed, so re-resolve the range) and apply the
# accent6 font color to the whole run.
$trEn2 = $shEn.TextFrame.TextRange
$foundEn2 = $trEn2.Find($newEn)
if ($foundEn2 -eq $null) {
  throw "Could not re-find the updated English comment #6 paragraph text"
}
$foundEn2.Font.Color.ObjectThemeColor = $accent6

# --- Part 2: Korean translation of comment #6 (shape "내용 개체 틀 11") ---
$shKo = $s.Shapes.Item(3)
$trKo = $shKo.TextFrame.TextRange

$runTexts = @(
  "6. ",
  "재순위",
  " 성능은 격자 ",
  "+ ",
  "트랜스포머",
  "((En)",
  "의 필요성과 트랜스포머를 사용하는 이유를 설명해야 합니다",
  ". ",
  "격자 ",
  "+ ",
  "트랜스포머",
  "(En)",
  "를 사용하여 성능을 향상시키는 방법이 설명되어야 합니다",
  ". ",
  "동일한 컴퓨팅 리소스에서 성능 비용은 무엇을 의미하나요",
  "? ",
  "격자 ",
  "+ ",
  "트랜스포머",
  "(En) ",
  "모델은 다이어그램을 통해 더 자세히 설명할 필요가 있습니다",
  "."
)

$cursor = 0
foreach ($txt in $runTexts) {
  $hit = $trKo.Find($txt, $cursor)
  if ($hit -eq $null) {
    throw "Could not find run text: $txt"
  }
  $hit.Font.Color.ObjectThemeColor = $accent6
  $cursor = $hit.Start + $hit.Length - 1
}
